# Apply permutation of species-record rows (2-27, row 19 unchanged) per commit diff.
# Each destination row receives the full data of a specific source row (matched by the
# "Id" identity that travels with the record), while constant columns (location, dates,
# observer, etc.) are untouched because they are identical for every row already.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= old row 18 (Id=111815492)
$ws.Range("A2").Value = 111815492
$ws.Range("B2").Value = 78578
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("I2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("Q2").Value = 457652
$ws.Range("R2").Value = 7058423
$ws.Range("AC2").Value = ""

# Row 3 <= old row 12 (Id=111815506)
$ws.Range("A3").Value = 111815506
$ws.Range("B3").Value = 77515
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("I3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = 457558
$ws.Range("R3").Value = 7057461
$ws.Range("AC3").Value = ""

# Row 4 <= old row 25 (Id=111815489)
$ws.Range("A4").Value = 111815489
$ws.Range("B4").Value = 56414
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("Q4").Value = 457851
$ws.Range("R4").Value = 7058248
$ws.Range("AC4").Value = "hack"

# Row 5 <= old row 14 (Id=111815486)
$ws.Range("A5").Value = 111815486
$ws.Range("B5").Value = 56398
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("Q5").Value = 457491
$ws.Range("R5").Value = 7057911
$ws.Range("AC5").Value = "ringhack"

# Row 6 <= old row 20 (Id=111815505)
$ws.Range("A6").Value = 111815505
$ws.Range("B6").Value = 77515
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("I6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("Q6").Value = 457628
$ws.Range("R6").Value = 7057503
$ws.Range("AC6").Value = ""

# Row 7 <= old row 13 (Id=111815498)
$ws.Range("A7").Value = 111815498
$ws.Range("B7").Value = 89423
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma"
$ws.Range("H7").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("Q7").Value = 457526
$ws.Range("R7").Value = 7057587
$ws.Range("AC7").Value = ""

# Row 8 <= old row 17 (Id=111815483)
$ws.Range("A8").Value = 111815483
$ws.Range("B8").Value = 56398
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("I8").Value = "'1"
$ws.Range("M8").Value = "födosökande"
$ws.Range("N8").Value = "observerad"
$ws.Range("Q8").Value = 457815
$ws.Range("R8").Value = 7058240
$ws.Range("AC8").Value = ""

# Row 9 <= old row 16 (Id=111815471)
$ws.Range("A9").Value = 111815471
$ws.Range("B9").Value = 90087
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 3298
$ws.Range("F9").Value = "Trådticka"
$ws.Range("G9").Value = "Climacocystis borealis"
$ws.Range("H9").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("Q9").Value = 457691
$ws.Range("R9").Value = 7058280
$ws.Range("AC9").Value = ""

# Row 10 <= old row 2 (Id=111815470)
$ws.Range("A10").Value = 111815470
$ws.Range("B10").Value = 90087
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 3298
$ws.Range("F10").Value = "Trådticka"
$ws.Range("G10").Value = "Climacocystis borealis"
$ws.Range("H10").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 457615
$ws.Range("R10").Value = 7058261
$ws.Range("AC10").Value = ""

# Row 11 <= old row 23 (Id=111815469)
$ws.Range("A11").Value = 111815469
$ws.Range("B11").Value = 90087
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 3298
$ws.Range("F11").Value = "Trådticka"
$ws.Range("G11").Value = "Climacocystis borealis"
$ws.Range("H11").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("Q11").Value = 457737
$ws.Range("R11").Value = 7057633
$ws.Range("AC11").Value = ""

# Row 12 <= old row 6 (Id=111815478)
$ws.Range("A12").Value = 111815478
$ws.Range("B12").Value = 90087
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 3298
$ws.Range("F12").Value = "Trådticka"
$ws.Range("G12").Value = "Climacocystis borealis"
$ws.Range("H12").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 457491
$ws.Range("R12").Value = 7057589
$ws.Range("AC12").Value = ""

# Row 13 <= old row 21 (Id=111815490)
$ws.Range("A13").Value = 111815490
$ws.Range("B13").Value = 56414
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 100049
$ws.Range("F13").Value = "Spillkråka"
$ws.Range("G13").Value = "Dryocopus martius"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("Q13").Value = 457487
$ws.Range("R13").Value = 7058060
$ws.Range("AC13").Value = "hack"

# Row 14 <= old row 27 (Id=111815501)
$ws.Range("A14").Value = 111815501
$ws.Range("B14").Value = 77515
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("I14").Value = ""
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("Q14").Value = 457502
$ws.Range("R14").Value = 7058397
$ws.Range("AC14").Value = ""

# Row 15 <= old row 22 (Id=111815500)
$ws.Range("A15").Value = 111815500
$ws.Range("B15").Value = 77515
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("I15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 457856
$ws.Range("R15").Value = 7058258
$ws.Range("AC15").Value = ""

# Row 16 <= old row 4 (Id=111815476)
$ws.Range("A16").Value = 111815476
$ws.Range("B16").Value = 90087
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 3298
$ws.Range("F16").Value = "Trådticka"
$ws.Range("G16").Value = "Climacocystis borealis"
$ws.Range("H16").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("Q16").Value = 457561
$ws.Range("R16").Value = 7058242
$ws.Range("AC16").Value = ""

# Row 17 <= old row 7 (Id=111815482)
$ws.Range("A17").Value = 111815482
$ws.Range("B17").Value = 56398
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("I17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("Q17").Value = 457734
$ws.Range("R17").Value = 7057882
$ws.Range("AC17").Value = "ringhack"

# Row 18 <= old row 15 (Id=111815495)
$ws.Range("A18").Value = 111815495
$ws.Range("B18").Value = 89423
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 5432
$ws.Range("F18").Value = "Granticka"
$ws.Range("G18").Value = "Porodaedalea chrysoloma"
$ws.Range("H18").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""
$ws.Range("Q18").Value = 457740
$ws.Range("R18").Value = 7057635
$ws.Range("AC18").Value = ""

# Row 20 <= old row 26 (Id=111815485)
$ws.Range("A20").Value = 111815485
$ws.Range("B20").Value = 56398
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("I20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("Q20").Value = 457447
$ws.Range("R20").Value = 7058136
$ws.Range("AC20").Value = "ringhack"

# Row 21 <= old row 10 (Id=111815494)
$ws.Range("A21").Value = 111815494
$ws.Range("B21").Value = 89419
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 1204
$ws.Range("F21").Value = "Gränsticka"
$ws.Range("G21").Value = "Phellopilus nigrolimitatus"
$ws.Range("H21").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("I21").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("Q21").Value = 457558
$ws.Range("R21").Value = 7057457
$ws.Range("AC21").Value = ""

# Row 22 <= old row 9 (Id=111815480)
$ws.Range("A22").Value = 111815480
$ws.Range("B22").Value = 90087
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 3298
$ws.Range("F22").Value = "Trådticka"
$ws.Range("G22").Value = "Climacocystis borealis"
$ws.Range("H22").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I22").Value = ""
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""
$ws.Range("Q22").Value = 457651
$ws.Range("R22").Value = 7057582
$ws.Range("AC22").Value = ""

# Row 23 <= old row 24 (Id=111815475)
$ws.Range("A23").Value = 111815475
$ws.Range("B23").Value = 90087
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 3298
$ws.Range("F23").Value = "Trådticka"
$ws.Range("G23").Value = "Climacocystis borealis"
$ws.Range("H23").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I23").Value = ""
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = ""
$ws.Range("Q23").Value = 457550
$ws.Range("R23").Value = 7058250
$ws.Range("AC23").Value = ""

# Row 24 <= old row 3 (Id=111815472)
$ws.Range("A24").Value = 111815472
$ws.Range("B24").Value = 90087
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 3298
$ws.Range("F24").Value = "Trådticka"
$ws.Range("G24").Value = "Climacocystis borealis"
$ws.Range("H24").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("I24").Value = ""
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = ""
$ws.Range("Q24").Value = 457859
$ws.Range("R24").Value = 7058252
$ws.Range("AC24").Value = ""

# Row 25 <= old row 11 (Id=111815503)
$ws.Range("A25").Value = 111815503
$ws.Range("B25").Value = 77515
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("I25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("Q25").Value = 457482
$ws.Range("R25").Value = 7057721
$ws.Range("AC25").Value = ""

# Row 26 <= old row 8 (Id=111815484)
$ws.Range("A26").Value = 111815484
$ws.Range("B26").Value = 56398
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 100109
$ws.Range("F26").Value = "Tretåig hackspett"
$ws.Range("G26").Value = "Picoides tridactylus"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("I26").Value = ""
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""
$ws.Range("Q26").Value = 457499
$ws.Range("R26").Value = 7058354
$ws.Range("AC26").Value = "ringhack"

# Row 27 <= old row 5 (Id=111815499)
$ws.Range("A27").Value = 111815499
$ws.Range("B27").Value = 89423
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = "Granticka"
$ws.Range("G27").Value = "Porodaedalea chrysoloma"
$ws.Range("H27").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
$ws.Range("Q27").Value = 457640
$ws.Range("R27").Value = 7057509
$ws.Range("AC27").Value = ""
